$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($r, $a, $b, $c, $d, $e, $f, $g, $h) {
    if ($a -ne $null) { $ws.Cells.Item($r, 1).Value = $a }
    if ($b -ne $null) { $ws.Cells.Item($r, 2).Value = $b }
    if ($c -ne $null) { $ws.Cells.Item($r, 3).Value = $c }
    if ($d -ne $null) { $ws.Cells.Item($r, 4).Value = $d }
    if ($e -ne $null) { $ws.Cells.Item($r, 5).Value = $e }
    if ($f -ne $null) { $ws.Cells.Item($r, 6).Value = $f }
    if ($g -ne $null) { $ws.Cells.Item($r, 7).Value = $g }
    if ($h -ne $null) { $ws.Cells.Item($r, 8).Value = $h }
}

# --- Update "Datos actualizados a ..." timestamp cell (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Octubre de 2020 a las 23:23"

# --- Update per-country case numbers (columns B-H); column A (country name)
#     unchanged unless noted. Rows are ordered by "Casos totales" (column B)
#     descending, so Botsuana's growth pushes it above "Republica de Yibuti"
#     and the labels for rows 127-131 shift down by one place. ---

# Row 4: Estados Unidos
Set-Row 4 $null 8438360 40790 5488734 2724529 0 323 225097

# Row 5: India
Set-Row 5 $null 7594284 46046 6730379 748670 0 593 115235

# Row 27: Israel
Set-Row 27 $null 304876 1767 278394 24219 0 54 2263

# Row 58: Barein
Set-Row 58 $null 78224 322 74683 3239 0 2 302

# Rows 127-131: Botsuana overtakes Republica de Yibuti, Nicaragua, Trinidad
# yTobago and Hong Kong in total-case ranking, so the country labels shift
# down one row.
Set-Row 127 "Botsuana" 5609 367 915 4673 0 1 21
Set-Row 128 "Republica de Yibuti" 5469 10 5379 29 0 0 61
Set-Row 129 "Nicaragua" 5353 0 4225 974 0 0 154
Set-Row 130 "Trinidad yTobago" 5298 1 3696 1505 0 1 97
Set-Row 131 "Hong Kong" 5257 15 4982 170 0 0 105

# Row 136: Ruanda
Set-Row 136 $null 4992 18 4797 161 0 0 34

# Row 159: Sierra Leona
Set-Row 159 $null 2331 1 1760 498 0 0 73

# Row 172: Curazao
Set-Row 172 $null 751 7 458 292 0 0 1
